$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.315.01"
$ws.Range("E2").Value = "  -0.16%  "

$ws.Range("D3").Value = "1.839.12"
$ws.Range("E3").Value = "  -0.33%  "

$ws.Range("D4").Value = "'0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.09%  "

$ws.Range("D5").Value = "'239.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "

$ws.Range("D6").Value = "'0.6275"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "'0.07430"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.83%  "

$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2891"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.35%  "

$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'24.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.38%  "

$ws.Range("D11").Value = "'0.07718"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.06%  "

$ws.Range("D12").Value = "1.836.79"
$ws.Range("E12").Value = "  -0.46%  "

$ws.Range("D13").Value = "'4.954"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.90%  "

$ws.Range("D14").Value = "'0.6739"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "

$ws.Range("D15").Value = "'0.00001019"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.00%  "

$ws.Range("D16").Value = "'81.67"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.54%  "

$ws.Range("D17").Value = "'6.215"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.92%  "

$ws.Range("D18").Value = "29.269.95"
$ws.Range("E18").Value = "  -0.42%  "

$ws.Range("D19").Value = "'229.21"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.05%  "

$ws.Range("E20").Value = "  -0.50%  "

$ws.Range("D21").Value = "'0.9998"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.12%  "

$ws.Range("D22").Value = "'7.351"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.49%  "

$ws.Range("E23").Value = "  +0.30%  "

$ws.Range("D24").Value = "'158.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.36%  "

$ws.Range("D25").Value = "'8.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.75%  "

$ws.Range("E26").Value = "  -2.15%  "

$ws.Range("E27").Value = "  -1.11%  "

$ws.Range("D28").Value = "'0.07312"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.97%  "

$ws.Range("D29").Value = "'1.454"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.95%  "

$ws.Range("D30").Value = "'1.475"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.44%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.036"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.39%  "

$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").Value = "'4.041"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.50%  "

$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("E34").Value = "  +0.03%  "

$ws.Range("D35").Value = "'0.6942"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.63%  "

$ws.Range("D36").Value = "'2.571"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.18%  "

$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.01835"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.43%  "

$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "'6.892"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.58%  "

$ws.Range("D39").Value = "'2.813"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.50%  "

$ws.Range("D40").Value = "1.236.47"
$ws.Range("E40").Value = "  -1.78%  "

$ws.Range("D41").Value = "'0.9344"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.27%  "

$ws.Range("D42").Value = "'0.9996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.12%  "

$ws.Range("D43").Value = "'100.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.39%  "

$ws.Range("D44").Value = "1.982.41"
$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("D45").Value = "'65.28"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("E46").Value = "  +4.68%  "

$ws.Range("D47").Value = "'1.701"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.79%  "

$ws.Range("E48").Value = "  -2.02%  "

$ws.Range("E49").Value = "  -2.96%  "

$ws.Range("D50").Value = "'8.856"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.55%  "

$ws.Range("D51").Value = "'0.3902"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.00%  "

